$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 7101.6665
$ws.Range("I7").Value = 1305
$ws.Range("J7").Value = 10000
$ws.Range("K7").Value = 1305
$ws.Range("L7").Value = 10000
$ws.Range("M7").Value = -1193
$ws.Range("N7").Value = -10224

$ws.Range("H14").Value = 7101.6665
$ws.Range("I14").Value = 1305
$ws.Range("J14").Value = 10000
$ws.Range("K14").Value = 1305
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = -1114
$ws.Range("N14").Value = -10382

$ws.Range("H34").Value = 1799.3334
$ws.Range("I34").Value = 1799.3334
$ws.Range("K34").Value = 1799.3334
$ws.Range("M34").Value = -1596.3334

$ws.Range("H36").Value = 1799.3334
$ws.Range("I36").Value = 1799.3334
$ws.Range("K36").Value = 1799.3334
$ws.Range("M36").Value = -1084.3334

$ws.Range("H70").Value = 5243.5
$ws.Range("I70").Value = 2491.6667
$ws.Range("J70").Value = 6160.778
$ws.Range("K70").Value = 7475.000100000001
$ws.Range("L70").Value = 18482.334
$ws.Range("M70").Value = -7205.000100000001
$ws.Range("N70").Value = -19022.334

$ws.Range("H73").Value = 5243.5
$ws.Range("I73").Value = 2491.6667
$ws.Range("J73").Value = 6160.778
$ws.Range("K73").Value = 7475.000100000001
$ws.Range("L73").Value = 18482.334
$ws.Range("M73").Value = -6539.000100000001
$ws.Range("N73").Value = -20354.334

$ws.Range("H132").Value = 7956.2905
$ws.Range("I132").Value = 7884.3105
$ws.Range("K132").Value = 23652.9315
$ws.Range("M132").Value = -21122.9315

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1073.08
$ws.Range("I32").Value = 992.4783
$ws.Range("K32").Value = 992.4783
$ws.Range("M32").Value = -705.4783

$ws.Range("H45").Value = 3867
$ws.Range("I45").Value = 2222
$ws.Range("J45").Value = 4525
$ws.Range("K45").Value = 2222
$ws.Range("L45").Value = 4525
$ws.Range("M45").Value = -1845
$ws.Range("N45").Value = -5279

$ws.Range("H61").Value = 1221.4
$ws.Range("I61").Value = 1221.4
$ws.Range("K61").Value = 1221.4
$ws.Range("M61").Value = -1009.4

$ws.Range("H123").Value = 1979899
$ws.Range("J123").Value = 1979899
$ws.Range("L123").Value = 1979899
$ws.Range("N123").Value = -1989699

$ws.Range("H132").Value = 1650
$ws.Range("I132").Value = 1650
$ws.Range("K132").Value = 4950
$ws.Range("M132").Value = -2420

$ws.Range("H136").Value = 1221.4
$ws.Range("I136").Value = 1221.4
$ws.Range("K136").Value = 3664.2
$ws.Range("M136").Value = -1114.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 5004.5
$ws.Range("I107").Value = 1674.1666
$ws.Range("K107").Value = 1674.1666
$ws.Range("M107").Value = 245.8334

$ws.Range("H126").Value = 50000
$ws.Range("J126").Value = 50000
$ws.Range("L126").Value = 50000
$ws.Range("N126").Value = -59880

$ws.Range("H134").Value = 1028.4166
$ws.Range("I134").Value = 1028.4166
$ws.Range("K134").Value = 3085.2498
$ws.Range("M134").Value = -550.2498000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6169.857
$ws.Range("I31").Value = 1854.2858
$ws.Range("K31").Value = 1854.2858
$ws.Range("M31").Value = -1559.2858

$ws.Range("H34").Value = 6169.857
$ws.Range("I34").Value = 1854.2858
$ws.Range("K34").Value = 1854.2858
$ws.Range("M34").Value = -1652.2858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 109.21429
$ws.Range("I2").Value = 108.416664
$ws.Range("J2").Value = 114
$ws.Range("K2").Value = 650.499984
$ws.Range("L2").Value = 684
$ws.Range("M2").Value = -537.499984
$ws.Range("N2").Value = -910

$ws.Range("H9").Value = 25150.5
$ws.Range("J9").Value = 300
$ws.Range("L9").Value = 900
$ws.Range("N9").Value = -1348

$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws.Range("H33").Value = 122.666664
$ws.Range("I33").Value = 139.5
$ws.Range("J33").Value = 89
$ws.Range("K33").Value = 837
$ws.Range("L33").Value = 534
$ws.Range("M33").Value = -554
$ws.Range("N33").Value = -1100

$ws.Range("H55").Value = 6529.909
$ws.Range("J55").Value = 6529.909
$ws.Range("L55").Value = 19589.727
$ws.Range("N55").Value = -19943.727

$ws.Range("H113").Value = 680.4
$ws.Range("I113").Value = 417.33334
$ws.Range("J113").Value = 1075
$ws.Range("K113").Value = 1252.00002
$ws.Range("L113").Value = 3225
$ws.Range("M113").Value = 917.9999800000001
$ws.Range("N113").Value = -7565

$ws.Range("H131").Value = 1227.5
$ws.Range("J131").Value = 1709.5
$ws.Range("L131").Value = 5128.5
$ws.Range("N131").Value = -15208.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 9937.200000000001
$ws.Range("J33").Value = 9937.200000000001
$ws.Range("L33").Value = 9937.200000000001
$ws.Range("N33").Value = -10441.2

$ws.Range("H126").Value = 2528.8462
$ws.Range("I126").Value = 1914.5834
$ws.Range("K126").Value = 5743.7502
$ws.Range("M126").Value = -3273.7502

$ws.Range("H132").Value = 201162
$ws.Range("I132").Value = 201162
$ws.Range("K132").Value = 603486
$ws.Range("M132").Value = -600956

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H136").Value = 2000
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 6000
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 45000
$ws.Range("J70").Value = 45000
$ws.Range("L70").Value = 45000
$ws.Range("N70").Value = -45630

$ws.Range("H73").Value = 45000
$ws.Range("J73").Value = 45000
$ws.Range("L73").Value = 45000
$ws.Range("N73").Value = -47184

$ws.Range("H136").Value = 2707
$ws.Range("I136").Value = 2267.4546
$ws.Range("K136").Value = 6802.3638
$ws.Range("M136").Value = -4252.3638
